$ws = $excel.ActiveWorkbook.ActiveSheet

$ws.Cells.Item(2, 3).Value = 1.285995721817017
$ws.Cells.Item(2, 5).Value = 633.1376521785569
$ws.Cells.Item(2, 6).Value = 0.02506652643680683
$ws.Cells.Item(2, 7).Value = 0.02096507363564691
$ws.Cells.Item(2, 8).Value = 0.01894712953950388
$ws.Cells.Item(2, 9).Value = 0.017012726749029
$ws.Cells.Item(2, 10).Value = 0.01590296042296959
$ws.Cells.Item(2, 11).Value = 0.0152148177820438
$ws.Cells.Item(2, 12).Value = 0.01449462711008997
$ws.Cells.Item(2, 13).Value = 0.01399981593082248
$ws.Cells.Item(2, 14).Value = 0.01373185237750325
$ws.Cells.Item(2, 15).Value = 0.0134696193681278
$ws.Cells.Item(2, 16).Value = 0.01324352616088779
$ws.Cells.Item(2, 17).Value = 0.01305243450386742
$ws.Cells.Item(2, 18).Value = 0.01280594221880037
$ws.Cells.Item(2, 19).Value = 0.01260862475379426
$ws.Cells.Item(2, 20).Value = 0.01251435662594541
$ws.Cells.Item(2, 21).Value = 0.01249285894215446
$ws.Cells.Item(2, 22).Value = 0.01242224824361597
$ws.Cells.Item(2, 23).Value = 0.01239743650215534
$ws.Cells.Item(2, 24).Value = 0.01236360435322458
$ws.Cells.Item(2, 25).Value = 0.01234186456488415

$ws.Cells.Item(3, 3).Value = 1.130040884017944
$ws.Cells.Item(3, 5).Value = 640.5298924388535
$ws.Cells.Item(3, 6).Value = 0.02572481502910171
$ws.Cells.Item(3, 7).Value = 0.02144421563181815
$ws.Cells.Item(3, 8).Value = 0.01883694509192812
$ws.Cells.Item(3, 9).Value = 0.01716556135325145
$ws.Cells.Item(3, 10).Value = 0.0161890877338635
$ws.Cells.Item(3, 11).Value = 0.01575714153878772
$ws.Cells.Item(3, 12).Value = 0.01457763671767395
$ws.Cells.Item(3, 13).Value = 0.01431542364540341
$ws.Cells.Item(3, 14).Value = 0.01362398399578339
$ws.Cells.Item(3, 15).Value = 0.01343654537485237
$ws.Cells.Item(3, 16).Value = 0.01319954878476851
$ws.Cells.Item(3, 17).Value = 0.01293279671630961
$ws.Cells.Item(3, 18).Value = 0.01283629422804617
$ws.Cells.Item(3, 19).Value = 0.01278314399732428
$ws.Cells.Item(3, 20).Value = 0.01265838314201972
$ws.Cells.Item(3, 21).Value = 0.01262527961413184
$ws.Cells.Item(3, 22).Value = 0.01261153688993961
$ws.Cells.Item(3, 23).Value = 0.01251064126113925
$ws.Cells.Item(3, 24).Value = 0.01251064126113925
$ws.Cells.Item(3, 25).Value = 0.01248596281557219

$ws.Cells.Item(4, 3).Value = 1.170999526977539
$ws.Cells.Item(4, 5).Value = 642.9889433988556
$ws.Cells.Item(4, 6).Value = 0.02533029965480071
$ws.Cells.Item(4, 7).Value = 0.021694161286074
$ws.Cells.Item(4, 8).Value = 0.01877384726451413
$ws.Cells.Item(4, 9).Value = 0.01739055671421734
$ws.Cells.Item(4, 10).Value = 0.01632338593675783
$ws.Cells.Item(4, 11).Value = 0.01540976419500897
$ws.Cells.Item(4, 12).Value = 0.01443983313907569
$ws.Cells.Item(4, 13).Value = 0.01416885987334103
$ws.Cells.Item(4, 14).Value = 0.01400713672301314
$ws.Cells.Item(4, 15).Value = 0.01357829700370496
$ws.Cells.Item(4, 16).Value = 0.0132926332438165
$ws.Cells.Item(4, 17).Value = 0.01319062221638451
$ws.Cells.Item(4, 18).Value = 0.01305811233564748
$ws.Cells.Item(4, 19).Value = 0.01290029029367325
$ws.Cells.Item(4, 20).Value = 0.0127757986721388
$ws.Cells.Item(4, 21).Value = 0.01272431987073791
$ws.Cells.Item(4, 22).Value = 0.01264178296633796
$ws.Cells.Item(4, 23).Value = 0.01257467422265225
$ws.Cells.Item(4, 24).Value = 0.01255676720952726
$ws.Cells.Item(4, 25).Value = 0.01253389753214143

$ws.Cells.Item(5, 3).Value = 1.107003450393677
$ws.Cells.Item(5, 5).Value = 645.0813347290641
$ws.Cells.Item(5, 6).Value = 0.02496373904142126
$ws.Cells.Item(5, 7).Value = 0.02020922609654791
$ws.Cells.Item(5, 8).Value = 0.01787491603163278
$ws.Cells.Item(5, 9).Value = 0.016604100222137
$ws.Cells.Item(5, 10).Value = 0.01598656969382061
$ws.Cells.Item(5, 11).Value = 0.01562740760328337
$ws.Cells.Item(5, 12).Value = 0.01503259423480453
$ws.Cells.Item(5, 13).Value = 0.01442933401418631
$ws.Cells.Item(5, 14).Value = 0.01389850777739996
$ws.Cells.Item(5, 15).Value = 0.01372355506486598
$ws.Cells.Item(5, 16).Value = 0.01337122787023335
$ws.Cells.Item(5, 17).Value = 0.01312490525672691
$ws.Cells.Item(5, 18).Value = 0.01303341716639262
$ws.Cells.Item(5, 19).Value = 0.01287796576222694
$ws.Cells.Item(5, 20).Value = 0.0127879904499463
$ws.Cells.Item(5, 21).Value = 0.01277476894600127
$ws.Cells.Item(5, 22).Value = 0.01268014742876622
$ws.Cells.Item(5, 23).Value = 0.01259609595031315
$ws.Cells.Item(5, 24).Value = 0.01258695232640611
$ws.Cells.Item(5, 25).Value = 0.01257468488750612

$ws.Cells.Item(6, 3).Value = 1.071038007736206
$ws.Cells.Item(6, 5).Value = 647.3317875380944
$ws.Cells.Item(6, 6).Value = 0.02496434351269935
$ws.Cells.Item(6, 7).Value = 0.02053236950047181
$ws.Cells.Item(6, 8).Value = 0.01846810989292908
$ws.Cells.Item(6, 9).Value = 0.01719210388035509
$ws.Cells.Item(6, 10).Value = 0.01634570265765975
$ws.Cells.Item(6, 11).Value = 0.01546236659711027
$ws.Cells.Item(6, 12).Value = 0.01473050443034872
$ws.Cells.Item(6, 13).Value = 0.01427442028942872
$ws.Cells.Item(6, 14).Value = 0.0138447443076013
$ws.Cells.Item(6, 15).Value = 0.0137206245310269
$ws.Cells.Item(6, 16).Value = 0.01343760974796686
$ws.Cells.Item(6, 17).Value = 0.01326762796419172
$ws.Cells.Item(6, 18).Value = 0.01312163039056797
$ws.Cells.Item(6, 19).Value = 0.0130312329651655
$ws.Cells.Item(6, 20).Value = 0.01290301585283488
$ws.Cells.Item(6, 21).Value = 0.01280316023512724
$ws.Cells.Item(6, 22).Value = 0.01272449425828077
$ws.Cells.Item(6, 23).Value = 0.01267624941514561
$ws.Cells.Item(6, 24).Value = 0.01265396120201057
$ws.Cells.Item(6, 25).Value = 0.01261855336331568

$ws.Cells.Item(7, 3).Value = 1.222987651824951
$ws.Cells.Item(7, 5).Value = 638.9283868123821
$ws.Cells.Item(7, 6).Value = 0.02577081560976348
$ws.Cells.Item(7, 7).Value = 0.0216657799316442
$ws.Cells.Item(7, 8).Value = 0.01960871139900268
$ws.Cells.Item(7, 9).Value = 0.01765868149068938
$ws.Cells.Item(7, 10).Value = 0.01646315123075331
$ws.Cells.Item(7, 11).Value = 0.01569204858315635
$ws.Cells.Item(7, 12).Value = 0.01497791783155951
$ws.Cells.Item(7, 13).Value = 0.01448573578478835
$ws.Cells.Item(7, 14).Value = 0.01418964511833517
$ws.Cells.Item(7, 15).Value = 0.01399874397796758
$ws.Cells.Item(7, 16).Value = 0.01365162263070663
$ws.Cells.Item(7, 17).Value = 0.01323305926687337
$ws.Cells.Item(7, 18).Value = 0.01303937315846117
$ws.Cells.Item(7, 19).Value = 0.01291799222295307
$ws.Cells.Item(7, 20).Value = 0.01272534919841801
$ws.Cells.Item(7, 21).Value = 0.01261868148759881
$ws.Cells.Item(7, 22).Value = 0.01258948923084972
$ws.Cells.Item(7, 23).Value = 0.01255882836664269
$ws.Cells.Item(7, 24).Value = 0.01249207355116681
$ws.Cells.Item(7, 25).Value = 0.01245474438230764

$ws.Cells.Item(8, 3).Value = 1.070997476577759
$ws.Cells.Item(8, 5).Value = 645.9737863574956
$ws.Cells.Item(8, 6).Value = 0.02509664537932519
$ws.Cells.Item(8, 7).Value = 0.02064058871609874
$ws.Cells.Item(8, 8).Value = 0.01879339777830611
$ws.Cells.Item(8, 9).Value = 0.01634441988545732
$ws.Cells.Item(8, 10).Value = 0.01593517398939657
$ws.Cells.Item(8, 11).Value = 0.01511178522747861
$ws.Cells.Item(8, 12).Value = 0.01479958644389083
$ws.Cells.Item(8, 13).Value = 0.01434843839819418
$ws.Cells.Item(8, 14).Value = 0.01402907806418205
$ws.Cells.Item(8, 15).Value = 0.01371995521328842
$ws.Cells.Item(8, 16).Value = 0.01337732136142214
$ws.Cells.Item(8, 17).Value = 0.01319026234793041
$ws.Cells.Item(8, 18).Value = 0.01308221559371382
$ws.Cells.Item(8, 19).Value = 0.01297225084760422
$ws.Cells.Item(8, 20).Value = 0.01288597888426684
$ws.Cells.Item(8, 21).Value = 0.01277112364609386
$ws.Cells.Item(8, 22).Value = 0.01267279132754145
$ws.Cells.Item(8, 23).Value = 0.0126461915766716
$ws.Cells.Item(8, 24).Value = 0.01260783033456349
$ws.Cells.Item(8, 25).Value = 0.01259208160540927

$ws.Cells.Item(9, 3).Value = 1.085999250411987
$ws.Cells.Item(9, 5).Value = 630.5421368071147
$ws.Cells.Item(9, 6).Value = 0.02536469901018662
$ws.Cells.Item(9, 7).Value = 0.02107642890649879
$ws.Cells.Item(9, 8).Value = 0.01888716515414414
$ws.Cells.Item(9, 9).Value = 0.01756657767343197
$ws.Cells.Item(9, 10).Value = 0.01646779357308875
$ws.Cells.Item(9, 11).Value = 0.01520907990740481
$ws.Cells.Item(9, 12).Value = 0.014747154611078
$ws.Cells.Item(9, 13).Value = 0.01403351134574128
$ws.Cells.Item(9, 14).Value = 0.01377001471849781
$ws.Cells.Item(9, 15).Value = 0.01320962905009281
$ws.Cells.Item(9, 16).Value = 0.01298904436730185
$ws.Cells.Item(9, 17).Value = 0.01280191597311973
$ws.Cells.Item(9, 18).Value = 0.0127686751323993
$ws.Cells.Item(9, 19).Value = 0.01262715605620978
$ws.Cells.Item(9, 20).Value = 0.01252553276166832
$ws.Cells.Item(9, 21).Value = 0.01241450248227665
$ws.Cells.Item(9, 22).Value = 0.01241316810278063
$ws.Cells.Item(9, 23).Value = 0.01233292028254955
$ws.Cells.Item(9, 24).Value = 0.0123140416216552
$ws.Cells.Item(9, 25).Value = 0.01229126972333557

$ws.Cells.Item(10, 3).Value = 1.098011493682861
$ws.Cells.Item(10, 5).Value = 640.1399326716473
$ws.Cells.Item(10, 6).Value = 0.02567667097771828
$ws.Cells.Item(10, 7).Value = 0.02142370810849646
$ws.Cells.Item(10, 8).Value = 0.01920397099097458
$ws.Cells.Item(10, 9).Value = 0.01747471344576185
$ws.Cells.Item(10, 10).Value = 0.01673297666007305
$ws.Cells.Item(10, 11).Value = 0.0158744398492185
$ws.Cells.Item(10, 12).Value = 0.01484234921993942
$ws.Cells.Item(10, 13).Value = 0.01436749102873098
$ws.Cells.Item(10, 14).Value = 0.01366197005633437
$ws.Cells.Item(10, 15).Value = 0.01331393723575329
$ws.Cells.Item(10, 16).Value = 0.01315861011905459
$ws.Cells.Item(10, 17).Value = 0.0129691829516327
$ws.Cells.Item(10, 18).Value = 0.01285979709735595
$ws.Cells.Item(10, 19).Value = 0.01275018371575741
$ws.Cells.Item(10, 20).Value = 0.01270095873092895
$ws.Cells.Item(10, 21).Value = 0.01257935930388868
$ws.Cells.Item(10, 22).Value = 0.01255026748680999
$ws.Cells.Item(10, 23).Value = 0.0125239321367929
$ws.Cells.Item(10, 24).Value = 0.01250508524413416
$ws.Cells.Item(10, 25).Value = 0.01247836126065589

$ws.Cells.Item(11, 3).Value = 1.03163480758667
$ws.Cells.Item(11, 5).Value = 636.0683877261199
$ws.Cells.Item(11, 6).Value = 0.02451370915961799
$ws.Cells.Item(11, 7).Value = 0.02025148068847606
$ws.Cells.Item(11, 8).Value = 0.01763290449258119
$ws.Cells.Item(11, 9).Value = 0.01663012414003038
$ws.Cells.Item(11, 10).Value = 0.0158508098660964
$ws.Cells.Item(11, 11).Value = 0.01515839230013833
$ws.Cells.Item(11, 12).Value = 0.01437238650879111
$ws.Cells.Item(11, 13).Value = 0.01412860805177591
$ws.Cells.Item(11, 14).Value = 0.01362686947749443
$ws.Cells.Item(11, 15).Value = 0.01344704872053002
$ws.Cells.Item(11, 16).Value = 0.01320620292463987
$ws.Cells.Item(11, 17).Value = 0.01303094587515947
$ws.Cells.Item(11, 18).Value = 0.01284420240339365
$ws.Cells.Item(11, 19).Value = 0.01277995326297816
$ws.Cells.Item(11, 20).Value = 0.01263443500611844
$ws.Cells.Item(11, 21).Value = 0.0125634456674783
$ws.Cells.Item(11, 22).Value = 0.01248393185608044
$ws.Cells.Item(11, 23).Value = 0.01242542362779569
$ws.Cells.Item(11, 24).Value = 0.01242542362779569
$ws.Cells.Item(11, 25).Value = 0.01239899391278986
